$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.982.73'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '3.788.89'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '705.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.46%  '
$ws.Range('D7').Value = '3.786.50'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.163'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.46'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000257'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('D15').Value = '4.427.84'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '3.788.88'
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '70.992.58'
$ws.Range('E17').Value = '  +2.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +17.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '482.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.717'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000146'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.57%  '
$ws.Range('D29').Value = '3.939.78'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +16.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.33'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.58'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.56'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.07%  '
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.21'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '3.739.35'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +13.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000328'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +23.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.967'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '160.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '49.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.40%  '
$ws.Range('E51').Value = '  +1.86%  '
